{"js": "// Site update for FAQ\n// Update the recorded simulation run timestamps/duration in the\n// \"Simulation parameters\" section of the document body.\n\nconst replacements = [\n  { find: \"Start time: 2017-12-27 18:32:32\", replace: \"Start time: 2018-01-31 12:35:14\" },\n  { find: \"End time: 2017-12-27 18:32:54\", replace: \"End time: 2018-01-31 12:35:37\" },\n  { find: \"Duration: 21.59 secs\", replace: \"Duration: 22.99 secs\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Site update for FAQ\n# Update the recorded simulation run timestamps/duration in the\n# \"Simulation parameters\" section of the document body.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstOccurrence($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1) | Out-Null\n}\n\nReplace-FirstOccurrence \"Start time: 2017-12-27 18:32:32\" \"Start time: 2018-01-31 12:35:14\"\nReplace-FirstOccurrence \"End time: 2017-12-27 18:32:54\" \"End time: 2018-01-31 12:35:37\"\nReplace-FirstOccurrence \"Duration: 21.59 secs\" \"Duration: 22.99 secs\"\n"}
